$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.182.78"
$ws.Range("E2").Value = "  +2.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.681.09"
$ws.Range("E3").Value = "  +7.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.54"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.25"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.670.72"
$ws.Range("E7").Value = "  +7.33%  "

$ws.Range("E8").Value = "  +4.11%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.200"
$ws.Range("E10").Value = "  -0.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.87"
$ws.Range("E11").Value = "  +25.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.611"
$ws.Range("E12").Value = "  +4.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.96"
$ws.Range("E13").Value = "  -0.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000290"
$ws.Range("E14").Value = "  +2.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.272.57"
$ws.Range("E15").Value = "  +7.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "675.11"
$ws.Range("E16").Value = "  -2.35%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "9.02"
$ws.Range("E17").Value = "  +4.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.687.22"
$ws.Range("E18").Value = "  +7.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.216.19"
$ws.Range("E19").Value = "  +2.07%  "

$ws.Range("E20").Value = "  +0.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.93"
$ws.Range("E21").Value = "  +1.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.56"
$ws.Range("E22").Value = "  +1.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.940"
$ws.Range("E23").Value = "  +4.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.28"
$ws.Range("E24").Value = "  +2.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.66"
$ws.Range("E25").Value = "  +0.91%  "

$ws.Range("E26").Value = "  +1.17%  "

$ws.Range("E27").Value = "  +6.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.12"
$ws.Range("E28").Value = "  +5.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.04"
$ws.Range("E30").Value = "  +4.61%  "

$ws.Range("E31").Value = "  +4.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.11"
$ws.Range("E32").Value = "  +4.00%  "

$ws.Range("E33").Value = "  -2.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.58"
$ws.Range("E34").Value = "  +6.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.01"
$ws.Range("E35").Value = "  +7.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "584.47"
$ws.Range("E36").Value = "  +1.84%  "

$ws.Range("E37").Value = "  +1.60%  "

$ws.Range("E38").Value = "  +4.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.72"
$ws.Range("E39").Value = "  +0.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("E41").Value = "  +9.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.620.43"
$ws.Range("E42").Value = "  +0.61%  "

$ws.Range("E43").Value = "  +2.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₃0760"
$ws.Range("E45").Value = "  +3.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "35.14"
$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("E47").Value = "  +2.83%  "

$ws.Range("E48").Value = "  +9.28%  "

$ws.Range("E49").Value = "  +4.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.71"
$ws.Range("E50").Value = "  +1.17%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.97"
$ws.Range("E51").Value = "  +9.11%  "
